# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" / "Valor Mora" block (rows 16-21, columns E/F) gets its
# row order reversed: what used to be listed newest-period-first (2311 down
# to 2306) is now listed oldest-period-first (2306 up to 2311), carrying the
# matching "Valor Mora" amount along with each period as the rows flip.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current (pre-edit) Periodo Mora / Valor Mora columns for the
# six data rows before overwriting anything. (NOTE: ".Value" getter is
# unreliable in this host - use ".Value2" for reads.)
$periodos = @()
$valores  = @()
for ($i = 16; $i -le 21; $i++) {
    $periodos += $ws.Range("E$i").Value2
    $valores  += $ws.Range("F$i").Value2
}

# Write them back in reverse row order, keeping Periodo <-> Valor paired.
$n = $periodos.Count
for ($i = 0; $i -lt $n; $i++) {
    $destRow = 16 + $i
    $srcIdx  = $n - 1 - $i
    $ws.Range("E$destRow").Value2 = [string]$periodos[$srcIdx]
    $ws.Range("F$destRow").Value2 = $valores[$srcIdx]
}
